$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 427, shifting existing rows 427:502 down to 428:503
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with the new data point
$ws.Cells.Item(427, 1).Value = 4
$ws.Cells.Item(427, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(427, 3).Value = "Los Lagos"
$ws.Cells.Item(427, 4).Value = 45258
$ws.Cells.Item(427, 5).Value = 10
$ws.Cells.Item(427, 6).Value = 100112032
$ws.Cells.Item(427, 7).Value = "Zapallo italiano"
$ws.Cells.Item(427, 8).Value = "Sin especificar"
$ws.Cells.Item(427, 9).Value = "Primera"
$ws.Cells.Item(427, 10).Value = 250
$ws.Cells.Item(427, 11).Value = 20000
$ws.Cells.Item(427, 12).Value = 20000
$ws.Cells.Item(427, 13).Value = 20000
$ws.Cells.Item(427, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(427, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(427, 16).Value = 400
$ws.Cells.Item(427, 17).Value = 50
$ws.Cells.Item(427, 18).Value = "Hortaliza"
